$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.740.10'
$ws.Range('E2').Value = '  -0.54%  '
$ws.Range('D3').Value = '2.282.32'
$ws.Range('E3').Value = '  -0.93%  '
$ws.Range('E4').Value = '  +0.07%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '114.65'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.86%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '266.56'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.34%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.643'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.49%  '
$ws.Range('E8').Value = '  +0.13%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.615'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -1.34%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '47.84'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -0.82%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0938'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.48%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '9.14'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.80%  '
$ws.Range('E13').Value = '  +1.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '15.45'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -3.23%  '
$ws.Range('D15').Value = '2.628.01'
$ws.Range('E15').Value = '  -0.78%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.877'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.35%  '
$ws.Range('D17').Value = '2.283.10'
$ws.Range('E17').Value = '  -0.74%  '
$ws.Range('D18').Value = '43.639.57'
$ws.Range('E18').Value = '  -0.53%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0000110'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.21%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.87'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.41%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '72.49'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.05%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '2.44'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.06'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.75%  '
$ws.Range('B24').Value = 'PancakeSwap'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.90'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.27%  '
$ws.Range('B25').Value = 'InternetComputer(DFINITY)'
$ws.Range('C25').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.52'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -2.01%  '
$ws.Range('E26').Value = '  +1.75%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.70'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.03%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '42.24'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.23%  '
$ws.Range('E29').Value = '  +0.63%  '
$ws.Range('E30').Value = '  -1.00%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '173.75'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -1.09%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '21.68'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +0.20%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0911'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -3.18%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.71'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -0.28%  '
$ws.Range('E35').Value = '  +1.71%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0386'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +6.05%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.68'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.14%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.94'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +3.66%  '
$ws.Range('E39').Value = '  -0.83%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.55'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +6.29%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '14.16'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +3.03%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '74.27'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.73%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.236'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.02%  '
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.13%  '
$ws.Range('B45').Value = 'THORChain'
$ws.Range('C45').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '5.96'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.37'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.86%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.28'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +3.16%  '
$ws.Range('B48').Value = 'FraxShare'
$ws.Range('C48').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '8.60'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -2.44%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.100'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.45%  '
$ws.Range('B50').Value = 'ordi'
$ws.Range('C50').Value = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '71.92'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +30.34%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '101.21'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.37%  '
